$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 5 (shifts old rows 5-19 down to 6-20)
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with data
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C5").Value = "Los Lagos"
$ws.Range("D5").Value = 44708
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 100112012
$ws.Range("G5").Value = "Espinaca"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 11000
$ws.Range("M5").Value = 11000
$ws.Range("N5").Value = "$/cuna 10 kilos"
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 1100
$ws.Range("Q5").Value = 10
$ws.Range("R5").Value = "Hortaliza"
